$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated simulation results for the 380 kV case
$ws.Range("B2").Value = 0.8933291878283285
$ws.Range("D2").Value = 0.003258175604951674
$ws.Range("E2").Value = 0.3575102998598396
$ws.Range("F2").Value = 0.5684327367567477
$ws.Range("G2").Value = 0.002362791442131427
$ws.Range("I2").Value = 0.6596824315284011
$ws.Range("L2").Value = 0.2991135008363699
$ws.Range("N2").Value = 1.362506658411533
$ws.Range("O2").Value = 1.815699715449085

$ws.Range("B3").Value = 0.8338038650850024
$ws.Range("D3").Value = 0.002973824538699432
$ws.Range("E3").Value = 0.3495491793288181
$ws.Range("F3").Value = 0.5432571050713761
$ws.Range("G3").Value = 0.002366191338576264
$ws.Range("I3").Value = 0.6718803821831688
$ws.Range("L3").Value = 0.2687331502148425
$ws.Range("N3").Value = 1.346805121431061
$ws.Range("O3").Value = 1.743156897276435

$ws.Range("B4").Value = 0.7974121339820499
$ws.Range("D4").Value = 0.002798406563439215
$ws.Range("E4").Value = 0.344681670660556
$ws.Range("F4").Value = 0.5282057394349664
$ws.Range("G4").Value = 0.002368392553594493
$ws.Range("I4").Value = 0.6798051901712682
$ws.Range("L4").Value = 0.2501100684635418
$ws.Range("N4").Value = 1.337671287969386
$ws.Range("O4").Value = 1.699953024009289

$ws.Range("B5").Value = 0.7826228868306089
$ws.Range("D5").Value = 0.002726719813987089
$ws.Range("E5").Value = 0.342703613017747
$ws.Range("F5").Value = 0.5221741110359801
$ws.Range("G5").Value = 0.002369318237647792
$ws.Range("I5").Value = 0.6831439808311204
$ws.Range("L5").Value = 0.2425290447543915
$ws.Range("N5").Value = 1.334077556814364
$ws.Range("O5").Value = 1.682682364823364

$ws.Range("B6").Value = 0.7801696354465832
$ws.Range("D6").Value = 0.002714804201504606
$ws.Range("E6").Value = 0.3423754985245964
$ws.Range("F6").Value = 0.5211787107052857
$ws.Range("G6").Value = 0.002369473681040755
$ws.Range("I6").Value = 0.6837049878778256
$ws.Range("L6").Value = 0.2412707169120551
$ws.Range("N6").Value = 1.333488600873494
$ws.Range("O6").Value = 1.67983480218038

$ws.Range("B7").Value = 0.7972125147449276
$ws.Range("D7").Value = 0.002797440584384248
$ws.Range("E7").Value = 0.344654971249696
$ws.Range("F7").Value = 0.5281239825070543
$ws.Range("G7").Value = 0.002368404921516668
$ws.Range("I7").Value = 0.6798497755419094
$ws.Range("L7").Value = 0.2500077951423947
$ws.Range("N7").Value = 1.337622300652598
$ws.Range("O7").Value = 1.699718749838809

$ws.Range("B8").Value = 0.8727729500101873
$ws.Range("D8").Value = 0.003160305559596566
$ws.Range("E8").Value = 0.3547611956327259
$ws.Range("F8").Value = 0.5596675983142347
$ws.Range("G8").Value = 0.002363940188982429
$ws.Range("I8").Value = 0.6637979138496867
$ws.Range("L8").Value = 0.2886322068806066
$ws.Range("N8").Value = 1.356988074152355
$ws.Range("O8").Value = 1.790408616347179

$ws.Range("B9").Value = 1.022148504779409
$ws.Range("D9").Value = 0.003865148055133716
$ws.Range("E9").Value = 0.3747312017924429
$ws.Range("F9").Value = 0.6247684301586389
$ws.Range("G9").Value = 0.002356082679966381
$ws.Range("I9").Value = 0.6357762637821343
$ws.Range("L9").Value = 0.364605843674326
$ws.Range("N9").Value = 1.398950876274995
$ws.Range("O9").Value = 1.978926762096137

$ws.Range("B10").Value = 1.132578952215226
$ws.Range("D10").Value = 0.004378689994993579
$ws.Range("E10").Value = 0.3894808207066518
$ws.Range("F10").Value = 0.6746059732758738
$ws.Range("G10").Value = 0.002350851423885198
$ws.Range("I10").Value = 0.6172987357268347
$ws.Range("L10").Value = 0.4205550790261725
$ws.Range("N10").Value = 1.432167830464067
$ws.Range("O10").Value = 2.124040777512278

$ws.Range("B11").Value = 1.182955657781235
$ws.Range("D11").Value = 0.004611339750130838
$ws.Range("E11").Value = 0.3962046144363924
$ws.Range("F11").Value = 0.6977212335151819
$ws.Range("G11").Value = 0.002348587997779721
$ws.Range("I11").Value = 0.6093516389756649
$ws.Range("L11").Value = 0.4460346869241789
$ws.Range("N11").Value = 1.44778860251688
$ws.Range("O11").Value = 2.191514966488853

$ws.Range("B12").Value = 1.202051283732715
$ws.Range("D12").Value = 0.004699295376568102
$ws.Range("E12").Value = 0.3987524847842394
$ws.Range("F12").Value = 0.7065386446418813
$ws.Range("G12").Value = 0.002347747528891025
$ws.Range("I12").Value = 0.6064082778525703
$ws.Range("L12").Value = 0.4556869110457455
$ws.Range("N12").Value = 1.453776338742557
$ws.Range("O12").Value = 2.217277269365013

$ws.Range("B13").Value = 1.197937871362797
$ws.Range("D13").Value = 0.00468035903776709
$ws.Range("E13").Value = 0.3982036840132679
$ws.Range("F13").Value = 0.704636798448476
$ws.Range("G13").Value = 0.002347927800153212
$ws.Range("I13").Value = 0.6070392447723502
$ws.Range("L13").Value = 0.4536079756005904
$ws.Range("N13").Value = 1.452483562553709
$ws.Range("O13").Value = 2.211719480321676

$ws.Range("B14").Value = 1.184526291637837
$ws.Range("D14").Value = 0.004618578824452158
$ws.Range("E14").Value = 0.3964141966778882
$ws.Range("F14").Value = 0.6984453594121192
$ws.Range("G14").Value = 0.002348518518776793
$ws.Range("I14").Value = 0.609108162868332
$ws.Range("L14").Value = 0.446828710240311
$ws.Range("N14").Value = 1.448279769398937
$ws.Range("O14").Value = 2.193630205245938

$ws.Range("B15").Value = 1.176313756652632
$ws.Range("D15").Value = 0.004580717777802334
$ws.Range("E15").Value = 0.39531829811893
$ws.Range("F15").Value = 0.6946612900152616
$ws.Range("G15").Value = 0.002348882516014096
$ws.Range("I15").Value = 0.6103840375535228
$ws.Range("L15").Value = 0.442676678490642
$ws.Range("N15").Value = 1.445714238960136
$ws.Range("O15").Value = 2.182577553626686

$ws.Range("B16").Value = 1.129289452050557
$ws.Range("D16").Value = 0.004363465953407086
$ws.Range("E16").Value = 0.3890416621477328
$ws.Range("F16").Value = 0.6731042975043522
$ws.Range("G16").Value = 0.002351001677390524
$ws.Range("I16").Value = 0.6178273313418995
$ws.Range("L16").Value = 0.4188904580894643
$ws.Range("N16").Value = 1.431157167706772
$ws.Range("O16").Value = 2.119660683150414

$ws.Range("B17").Value = 1.100476864659583
$ws.Range("D17").Value = 0.004229938548498069
$ws.Range("E17").Value = 0.3851945410606987
$ws.Range("F17").Value = 0.6599937121830379
$ws.Range("G17").Value = 0.002352331442459316
$ws.Range("I17").Value = 0.6225110432834264
$ws.Range("L17").Value = 0.404305288522977
$ws.Range("N17").Value = 1.422356925404102
$ws.Range("O17").Value = 2.081438288427364

$ws.Range("B18").Value = 1.083918005436999
$ws.Range("D18").Value = 0.004153046783663683
$ws.Range("E18").Value = 0.3829831263653105
$ws.Range("F18").Value = 0.6524946136892424
$ws.Range("G18").Value = 0.002353107239787286
$ws.Range("I18").Value = 0.6252481328301858
$ws.Range("L18").Value = 0.3959189447804476
$ws.Range("N18").Value = 1.417343334426306
$ws.Range("O18").Value = 2.059591232195316

$ws.Range("B19").Value = 1.078313796654754
$ws.Range("D19").Value = 0.004126997201808535
$ws.Range("E19").Value = 0.3822346204013627
$ws.Range("F19").Value = 0.6499627106418302
$ws.Range("G19").Value = 0.002353371794791323
$ws.Range("I19").Value = 0.6261822706724077
$ws.Range("L19").Value = 0.3930799413309956
$ws.Range("N19").Value = 1.415654102640332
$ws.Range("O19").Value = 2.052217766026388

$ws.Range("B20").Value = 1.103542638264855
$ws.Range("D20").Value = 0.004244162156837206
$ws.Range("E20").Value = 0.385603936480571
$ws.Range("F20").Value = 0.6613850316226717
$ws.Range("G20").Value = 0.002352188753948284
$ws.Range("I20").Value = 0.6220079882780825
$ws.Range("L20").Value = 0.4058576321725411
$ws.Range("N20").Value = 1.423288756967935
$ws.Range("O20").Value = 2.085492892980028

$ws.Range("B21").Value = 1.188465090450165
$ws.Range("D21").Value = 0.00463672911521229
$ws.Range("E21").Value = 0.396939768426833
$ws.Range("F21").Value = 0.7002621911298661
$ws.Range("G21").Value = 0.002348344559449878
$ws.Range("I21").Value = 0.6084986784967831
$ws.Range("L21").Value = 0.4488198487811417
$ws.Range("N21").Value = 1.449512564985099
$ws.Range("O21").Value = 2.198937719961918

$ws.Range("B22").Value = 1.24407733981127
$ws.Range("D22").Value = 0.004892453791033091
$ws.Range("E22").Value = 0.4043582445300089
$ws.Range("F22").Value = 0.7260448101941961
$ws.Range("G22").Value = 0.002345929118368359
$ws.Range("I22").Value = 0.6000545000748431
$ws.Range("L22").Value = 0.4769192379187928
$ws.Range("N22").Value = 1.467073458369271
$ws.Range("O22").Value = 2.274312593796139

$ws.Range("B23").Value = 1.214386339884697
$ws.Range("D23").Value = 0.004756047365606975
$ws.Range("E23").Value = 0.4003980661464723
$ws.Range("F23").Value = 0.7122497978454021
$ws.Range("G23").Value = 0.002347209439674291
$ws.Range("I23").Value = 0.6045260585002199
$ws.Range("L23").Value = 0.4619202599506309
$ws.Range("N23").Value = 1.457662543336141
$ws.Range("O23").Value = 2.233970473640568

$ws.Range("B24").Value = 1.102156583598344
$ws.Range("D24").Value = 0.004237732052761345
$ws.Range("E24").Value = 0.3854188477116338
$ws.Range("F24").Value = 0.6607558966652647
$ws.Range("G24").Value = 0.002352253228298657
$ws.Range("I24").Value = 0.6222352812718877
$ws.Range("L24").Value = 0.4051558211180009
$ws.Range("N24").Value = 1.422867333231522
$ws.Range("O24").Value = 2.083659409168263

$ws.Range("B25").Value = 0.98161482115475
$ws.Range("D25").Value = 0.003675211344010165
$ws.Range("E25").Value = 0.3693140804043082
$ws.Range("F25").Value = 0.6068065047297324
$ws.Range("G25").Value = 0.002358112816098878
$ws.Range("I25").Value = 0.6429865542096203
$ws.Range("L25").Value = 0.3440291964822961
$ws.Range("N25").Value = 1.387176561929522
$ws.Range("O25").Value = 1.926774207316839
